# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The workbook's "K" column (column G, header "K") previously held a
# Strike# derived count. The source data was regenerated upstream and the
# newly computed strikeout ("K") values are written back into column G
# for every game row (rows 2-58).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_vals for column G ("K"), keyed by worksheet row number.
$sVals = [ordered]@{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 2
    18 = 0
    19 = 3
    20 = 0
    21 = 0
    22 = 1
    23 = 2
    24 = 1
    25 = 1
    26 = 0
    27 = 1
    28 = 2
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 1
    34 = 0
    35 = 0
    36 = 2
    37 = 1
    38 = 1
    39 = 0
    40 = 1
    41 = 0
    42 = 2
    43 = 1
    44 = 1
    45 = 0
    46 = 1
    47 = 0
    48 = 0
    49 = 1
    50 = 3
    51 = 2
    52 = 1
    53 = 1
    54 = 0
    55 = 0
    56 = 2
    57 = 1
    58 = 1
}

foreach ($row in $sVals.Keys) {
    $ws.Cells.Item($row, 7).Value = $sVals[$row]
}
